$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (and values) from column R (rows 2-5) into new column S
$ws.Range("R2:R5").Copy($ws.Range("S2:S5"))

# Set the new column's values for the 2022 data point
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 211650
$ws.Range("S5").Value = 2.9794303052841493

# Update the active selection to match the new "last populated" cell in row 2
$ws.Range("S2").Select()
